# The sheet originally held a single numeric vector spread over rows 1-3
# and 5-7 (row 4 left blank as a spacer). The commit turns this into a
# "code chunk" style output: the R source for the chunk goes in row 1
# (in a monospace font), a blank spacer row follows, and the original
# vector output (now split into two sequential vectors) is pushed down
# two rows, keeping the same internal blank-row spacing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting two blank rows above row 1 shifts every existing row (and its
# contents/types) down by two, which reproduces the target layout exactly:
#   old row 1 (1)   -> row 3
#   old row 2 (2)   -> row 4
#   old row 3 (3)   -> row 5
#   old row 4 (gap) -> row 6
#   old row 5 (4)   -> row 7
#   old row 6 (5)   -> row 8
#   old row 7 (6)   -> row 9
$ws.Rows("1:2").Insert()

# Row 1 now becomes the source-code cell for the chunk; row 2 stays blank
# as the separator between the code echo and its printed output.
$ws.Range("A1").Value = "setNames(1:3, 1:3); setNames(4:6, 4:6)"
$ws.Range("A1").Font.Name = "Courier New"
